# Apply updated dSF (column F) values to the kennedy_ian.xlsx workbook.
# These reflect a "repull" of the underlying data where the final score
# differential (dSF) diverges from the starting score differential (dS0)
# for certain games.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    4  = -6
    6  = -5
    7  = -8
    8  = 5
    10 = 0
    13 = 2
    19 = -1
    25 = 0
    28 = -4
    32 = -2
    34 = 1
    39 = -1
    43 = -1
    51 = 0
    54 = 1
    55 = 1
    58 = -2
    59 = -8
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 6).Value = $updates[$row]
}
